$strXie = "Xie F, Pullenayegum E, Gaebel K, Bansback N, Bryan S, Ohinmaa A, et al. A Time Trade-off-derived Value Set of the EQ-5D-5L for Canada. Med Care [Internet]. 2016 Jan;54(1):98–105. Available from: http://dx.doi.org/10.1097/MLR.0000000000000447"
$strOpenManuscript = "Open access manuscript"
$strJanssen = "Janssen, B., & Szende, A. (2014). Population norms for the EQ‐5D. Chapter 3. In A. Szende, B. Janssen, & J. Cabases (Eds.), Self‐reported`npopulation health: An international perspective based on EQ‐5D. Amsterdam, Netherlands: Springer. 978‐94‐007‐7595‐4."
$strOpenReport = "Open access report"
$strNorway = "Norwegian Medicines Agency. (2018). Guidelines for the submission of documentation for single technology assessment (STA) of pharmaceuticals`n[Internet]. Available from https://legemiddelverket.no/Documents/English/Public%20funding%20and%20pricing/`nDocumentation%20for%20STA/Guidelines_april_2018.pdf"
$strToBeAdded2 = "To be added - contact manuscript author"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (Israel / age_covid row): "To be added" -> "To be added - contact manuscript author"
# Style changes from red-wrap (s=7) to plain-wrap (s=1); copy format from a donor cell that already uses style 1.
$ws.Range("A11").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D11").Value = $strToBeAdded2
$ws.Range("E11").Value = $strToBeAdded2

# --- Row 12 (UK, qol_norm row) ---
$ws.Range("E2").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D12").Value = $strJanssen
$ws.Range("E12").Value = $strOpenReport
$ws.Rows.Item(12).RowHeight = 145

# --- Row 13 (US, qol_norm row) ---
$ws.Range("E2").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D13").Value = $strJanssen
$ws.Range("E13").Value = $strOpenReport
$ws.Rows.Item(13).RowHeight = 145

# --- Row 14 (Canada, qol_norm row) ---
$ws.Range("E2").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D14").Value = $strXie
$ws.Range("E14").Value = $strOpenManuscript
$ws.Rows.Item(14).RowHeight = 145

# --- Row 15 (Norway, qol_norm row) ---
$ws.Range("E2").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D15").Value = $strNorway
$ws.Range("E15").Value = $strOpenReport
$ws.Rows.Item(15).RowHeight = 174

# --- Row 16 (Israel, qol_norm row): back to plain style (s=1), "To be added - contact manuscript author" ---
$ws.Range("A11").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D16").Value = $strToBeAdded2
$ws.Range("E16").Value = $strToBeAdded2

# --- Update selected cell shown when the workbook is reopened ---
$ws.Range("A3").Select() | Out-Null
